$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AK_1.png", 0.344, 0.001, 0.095, 0.367, "Tidak Diketahui", "Salah"),
    @("AK_2.png", 0.6929999999999999, 0.002, 0.337, 0.7, "Akhlak Kamiswara", "Benar"),
    @("AK_3.png", 0.672, 0.002, 0.427, 0.733, "Akhlak Kamiswara", "Benar"),
    @("AK_4.png", 0.335, 0.001, 0.43, 0.267, "Tidak Diketahui", "Salah"),
    @("AK_5.png", 0.335, 0.001, 0.439, 0.333, "Tidak Diketahui", "Salah"),
    @("MIB_1.png", 1.035, 0.003, 0.08500000000000001, 0.367, "Tidak Diketahui", "Salah"),
    @("MIB_2.png", 0.9320000000000001, 0.003, 0.275, 0.633, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_3.png", 1.182, 0.004, 0.546, 0.9330000000000001, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_4.png", 0.915, 0.003, 0.199, 0.733, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_5.png", 1.06, 0.004, 0.236, 0.767, "Muhammad Iqbal Baqi", "Benar"),
    @("AAH_1.png", 0.649, 0.002, 0.241, 0.9, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_2.png", 0.923, 0.003, 0.604, 1, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_3.png", 0.667, 0.002, 0.252, 0.967, "Andrea Ayunove Hutami", "Benar"),
    @("TI_1.png", 0.778, 0.003, 0.197, 0.533, "Muhammad Iqbal Baqi", "Salah"),
    @("TI_2.png", 0.82, 0.003, 0.327, 0.533, "Toni Ismail", "Benar"),
    @("TI_3.png", 0.5610000000000001, 0.002, 0.554, 0.433, "Tidak Diketahui", "Salah"),
    @("TI_4.png", 0.53, 0.002, 0.281, 0.6, "Toni Ismail", "Benar"),
    @("TI_5.png", 0.891, 0.003, 0.357, 0.533, "Toni Ismail", "Benar"),
    @("RAS_1.png", 0.457, 0.002, 0.371, 0.333, "Tidak Diketahui", "Salah"),
    @("RAS_2.png", 0.777, 0.003, 0.272, 0.5, "Ridha Ayu Salsabila", "Benar"),
    @("RAS_3.png", 0.358, 0.001, 0.343, 0.233, "Tidak Diketahui", "Salah"),
    @("RAS_4.png", 0.9409999999999999, 0.003, 0.099, 0.3, "Tidak Diketahui", "Salah"),
    @("RAS_5.png", 0.906, 0.003, 0.36, 0.533, "Ridha Ayu Salsabila", "Benar"),
    @("RR_1.png", 1.052, 0.003, 0.438, 0.667, "Rafiqo Rapitasari", "Benar"),
    @("RR_2.png", 1.057, 0.004, 0.49, 0.5669999999999999, "Rafiqo Rapitasari", "Benar"),
    @("RR_3.png", 0.748, 0.002, 0.08799999999999999, 0.5, "Rafiqo Rapitasari", "Benar"),
    @("RR_4.png", 1.033, 0.003, 0.443, 0.633, "Rafiqo Rapitasari", "Benar"),
    @("RR_5.png", 1.048, 0.003, 0.465, 0.733, "Rafiqo Rapitasari", "Benar"),
    @("AR_1.png", 0.677, 0.002, 0.317, 0.667, "Arizli Romadhon", "Benar"),
    @("GA_1.png", 1.03, 0.003, 0.441, 0.9, "Gege Ardiyansyah", "Benar"),
    @("GA_2.png", 0.495, 0.002, 0.25, 0.7, "Gege Ardiyansyah", "Benar"),
    @("GA_3.png", 0.5629999999999999, 0.002, 0.169, 0.6, "Gege Ardiyansyah", "Benar"),
    @("FY_1.png", 0.8159999999999999, 0.003, 0.239, 0.6, "Fanny Yusuf", "Benar"),
    @("FY_2.png", 1.224, 0.004, 0.497, 0.3, "Tidak Diketahui", "Salah"),
    @("FY_3.png", 1.122, 0.004, 0.305, 0.533, "Gege Ardiyansyah", "Salah"),
    @("FY_4.png", 0.985, 0.003, 0.23, 0.4, "Tidak Diketahui", "Salah"),
    @("TO_1.png", 0.629, 0.002, 0.357, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_2.png", 0.77, 0.003, 0.385, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_3.png", 0.655, 0.002, 0.391, 0.5, "Tiara Oktavian", "Benar"),
    @("TO_4.png", 1.841, 0.006, 0.08500000000000001, 0.5669999999999999, "Tidak Diketahui", "Salah"),
    @("TO_5.png", 1.802, 0.006, 0.076, 0.533, "Tiara Oktavian", "Benar"),
    @("TD_1.png", 1.375, 0.004, 0.024, 0.267, "Tidak Diketahui", "Benar"),
    @("TD_2.png", 1.34, 0.004, 0.04, 0.233, "Tidak Diketahui", "Benar"),
    @("TD_3.png", 0.896, 0.003, 0.174, 0.367, "Tidak Diketahui", "Benar"),
    @("TD_4.png", 0.761, 0.003, 0.099, 0.2, "Tidak Diketahui", "Benar")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
